$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format risky numeric-looking text cells as Text so Excel does not
# coerce strings like "544.93" or "0.0940" into numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "60.221.50"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").Value = "2.540.08"
$ws.Range("E3").Value = "  +2.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "544.93"
$ws.Range("E5").Value = "  +1.17%  "

# Row 6
$ws.Range("D6").Value = "145.34"
$ws.Range("E6").Value = "  -0.32%  "

# Row 7
$ws.Range("E7").Value = "  -0.43%  "

# Row 8
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  +0.24%  "

# Row 9
$ws.Range("D9").Value = "2.571.13"
$ws.Range("E9").Value = "  +3.53%  "

# Row 10
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +1.26%  "

# Row 11
$ws.Range("E11").Value = "  +1.53%  "

# Row 12
$ws.Range("D12").Value = "5.56"
$ws.Range("E12").Value = "  +4.32%  "

# Row 13
$ws.Range("D13").Value = "0.364"
$ws.Range("E13").Value = "  +2.13%  "

# Row 14
$ws.Range("D14").Value = "2.984.77"
$ws.Range("E14").Value = "  +2.59%  "

# Row 15
$ws.Range("D15").Value = "24.18"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
$ws.Range("D16").Value = "60.145.27"
$ws.Range("E16").Value = "  +1.74%  "

# Row 17
$ws.Range("E17").Value = "  +3.47%  "

# Row 18
$ws.Range("D18").Value = "2.580.34"
$ws.Range("E18").Value = "  +3.33%  "

# Row 19
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  -0.23%  "

# Row 20
$ws.Range("D20").Value = "4.36"
$ws.Range("E20").Value = "  +0.20%  "

# Row 21
$ws.Range("D21").Value = "329.19"
$ws.Range("E21").Value = "  +1.37%  "

# Row 22
$ws.Range("E22").Value = "  +0.47%  "

# Row 23
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +3.53%  "

# Row 24
$ws.Range("D24").Value = "62.56"
$ws.Range("E24").Value = "  +2.65%  "

# Row 25
$ws.Range("D25").Value = "0.442"
$ws.Range("E25").Value = "  -0.22%  "

# Row 26
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +3.87%  "

# Row 27
$ws.Range("D27").Value = "0.990"
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +3.08%  "

# Row 29
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  +0.55%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0803"
$ws.Range("E30").Value = "  +2.75%  "

# Row 31
$ws.Range("E31").Value = "  +0.10%  "

# Row 32
$ws.Range("E32").Value = "  -3.09%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.51"
$ws.Range("E33").Value = "  +6.04%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "162.16"
$ws.Range("E34").Value = "  +1.91%  "

# Row 35
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").Value = "18.82"
$ws.Range("E36").Value = "  +1.14%  "

# Row 37
$ws.Range("D37").Value = "4.50"
$ws.Range("E37").Value = "  +1.21%  "

# Row 38
$ws.Range("D38").Value = "1.66"
$ws.Range("E38").Value = "  +0.24%  "

# Row 39
$ws.Range("D39").Value = "5.81"
$ws.Range("E39").Value = "  -2.37%  "

# Row 40
$ws.Range("D40").Value = "37.26"
$ws.Range("E40").Value = "  +2.04%  "

# Row 41
$ws.Range("D41").Value = "304.86"
$ws.Range("E41").Value = "  -2.80%  "

# Row 42
$ws.Range("D42").Value = "0.843"
$ws.Range("E42").Value = "  +1.39%  "

# Row 43
$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  +0.29%  "

# Row 44
$ws.Range("E44").Value = "  -0.10%  "

# Row 45
$ws.Range("D45").Value = "0.609"
$ws.Range("E45").Value = "  +2.72%  "

# Row 46
$ws.Range("E46").Value = "  +0.84%  "

# Row 47
$ws.Range("D47").Value = "19.19"
$ws.Range("E47").Value = "  +4.40%  "

# Row 48
$ws.Range("D48").Value = "0.0940"
$ws.Range("E48").Value = "  +0.70%  "

# Row 49
$ws.Range("D49").Value = "124.77"
$ws.Range("E49").Value = "  -0.12%  "

# Row 50
$ws.Range("D50").Value = "0.0524"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
$ws.Range("D51").Value = "0.0230"
$ws.Range("E51").Value = "  +0.23%  "

# Restore the original (default/Normal) style on those cells so only the
# underlying text content changed, matching the source edit which left
# cell formatting untouched.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
